# Applies the "fix required fields on RC-DE" edit: updates the
# human-readable "Champ correspondant" labels in the schema tables.
#
# Mapping (Table, Row, Column -> new text):
#   T1 R2 C2 : messageId                                       -> Identifiant du message
#   T1 R5 C2 : Type de message NexSIS                          -> Type de message
#   T1 R6 C2 : Status du message                                -> Statut du message
#   T2 R2 C2 : Identifiant technique du système emetteur        -> Nom du système emetteur
#   T2 R3 C2 : uri                                               -> URI (identifiant technique)
#   T4 R3 C2 : uri                                               -> URI (identifiant technique)

$d = $word.ActiveDocument

function Set-CellText {
    param([int]$TableIndex, [int]$RowIndex, [int]$ColIndex, [string]$OldText, [string]$NewText)

    $table = $d.Tables.Item($TableIndex)
    $cell = $table.Cell($RowIndex, $ColIndex)
    $current = $cell.Range.Text -replace "[\x07\x0d]", ""

    if ($current -ne $OldText) {
        throw "Unexpected text in Table $TableIndex Row $RowIndex Col $ColIndex : '$current' (expected '$OldText')"
    }

    $cell.Range.Text = $NewText
}

Set-CellText 1 2 2 "messageId" "Identifiant du message"
Set-CellText 1 5 2 "Type de message NexSIS" "Type de message"
Set-CellText 1 6 2 "Status du message" "Statut du message"
Set-CellText 2 2 2 "Identifiant technique du système emetteur" "Nom du système emetteur"
Set-CellText 2 3 2 "uri" "URI (identifiant technique)"
Set-CellText 4 3 2 "uri" "URI (identifiant technique)"

Write-Output "Done."
